$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.481.69"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "1.942.42"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'243.17"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").Value = "'0.611"
$ws.Range("E6").Value = "  -0.40%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'57.10"
$ws.Range("E8").Value = "  -1.67%  "
$ws.Range("D9").Value = "'0.363"
$ws.Range("E9").Value = "  -1.90%  "
$ws.Range("D10").Value = "'0.0806"
$ws.Range("E10").Value = "  -2.23%  "
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("D12").Value = "'21.70"
$ws.Range("E12").Value = "  +0.62%  "
$ws.Range("D13").Value = "2.227.32"
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("D14").Value = "'0.805"
$ws.Range("E14").Value = "  -3.20%  "
$ws.Range("D15").Value = "'13.30"
$ws.Range("E15").Value = "  -1.94%  "
$ws.Range("D16").Value = "'5.17"
$ws.Range("E16").Value = "  -2.21%  "
$ws.Range("D17").Value = "1.942.31"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "36.462.01"
$ws.Range("E18").Value = "  +0.49%  "
$ws.Range("E19").Value = "  -1.06%  "
$ws.Range("D20").Value = "0.0₃0853"
$ws.Range("E20").Value = "  -2.45%  "
$ws.Range("D21").Value = "'227.41"
$ws.Range("E21").Value = "  -0.90%  "
$ws.Range("D22").Value = "'4.95"
$ws.Range("E22").Value = "  -1.30%  "
$ws.Range("E24").Value = "  -5.04%  "
$ws.Range("E25").Value = "  +0.77%  "
$ws.Range("D26").Value = "'9.16"
$ws.Range("E26").Value = "  -3.30%  "
$ws.Range("D27").Value = "'159.96"
$ws.Range("E27").Value = "  -2.33%  "
$ws.Range("E28").Value = "  +13.82%  "
$ws.Range("E29").Value = "  -2.05%  "
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "'1.09"
$ws.Range("E31").Value = "  -5.30%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'4.62"
$ws.Range("E32").Value = "  -2.21%  "
$ws.Range("D33").Value = "'0.0615"
$ws.Range("E33").Value = "  -3.18%  "
$ws.Range("D34").Value = "'4.15"
$ws.Range("E34").Value = "  -3.79%  "
$ws.Range("D35").Value = "'6.12"
$ws.Range("E35").Value = "  +0.58%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("E37").Value = "  -0.97%  "
$ws.Range("E38").Value = "  +1.73%  "
$ws.Range("D39").Value = "'3.27"
$ws.Range("E39").Value = "  +13.36%  "
$ws.Range("D40").Value = "'0.0986"
$ws.Range("E40").Value = "  -0.18%  "
$ws.Range("D41").Value = "'2.90"
$ws.Range("E41").Value = "  +1.04%  "
$ws.Range("D42").Value = "'0.0208"
$ws.Range("E42").Value = "  -0.94%  "
$ws.Range("D43").Value = "'1.14"
$ws.Range("E43").Value = "  -3.10%  "
$ws.Range("E44").Value = "  +0.38%  "
$ws.Range("D45").Value = "1.341.87"
$ws.Range("E45").Value = "  -0.16%  "
$ws.Range("D46").Value = "'1.02"
$ws.Range("E46").Value = "  -1.37%  "
$ws.Range("D47").Value = "'86.12"
$ws.Range("E47").Value = "  -2.55%  "
$ws.Range("E48").Value = "  -3.66%  "
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("D50").Value = "2.120.09"
$ws.Range("E50").Value = "  -0.50%  "
$ws.Range("D51").Value = "'42.99"
$ws.Range("E51").Value = "  -7.25%  "
